$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20 (shifts existing rows 20-43 down to 21-44)
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new weekly record
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C20").Value = "Los Lagos"
$ws.Range("D20").Value = Get-Date -Year 2021 -Month 12 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100103
$ws.Range("H20").Value = "Frutos de hueso (carozo)"
$ws.Range("I20").Value = 100103001
$ws.Range("J20").Value = "Cereza"
$ws.Range("K20").Value = "Santina"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 11000
$ws.Range("O20").Value = 12000
$ws.Range("P20").Value = 11500
$ws.Range("Q20").Value = "`$/caja 8 kilos"
$ws.Range("R20").Value = "Provincia de Curicó"
$ws.Range("S20").Value = 1438
$ws.Range("T20").Value = 8
